$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 9).Value = "sv"
$ws.Cells.Item(7, 10).Value = "Statement-opinion"
$ws.Cells.Item(9, 9).Value = "sd"
$ws.Cells.Item(9, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(14, 9).Value = "aa"
$ws.Cells.Item(14, 10).Value = "Agree/Accept"
$ws.Cells.Item(16, 9).Value = "aa"
$ws.Cells.Item(16, 10).Value = "Agree/Accept"
$ws.Cells.Item(19, 9).Value = "sd"
$ws.Cells.Item(19, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(20, 9).Value = "sd"
$ws.Cells.Item(20, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(21, 9).Value = "sd"
$ws.Cells.Item(21, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(24, 9).Value = "ba"
$ws.Cells.Item(24, 10).Value = "Appreciation"
$ws.Cells.Item(25, 9).Value = "sd"
$ws.Cells.Item(25, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(28, 9).Value = "%"
$ws.Cells.Item(28, 10).Value = "Uninterpretable"
$ws.Cells.Item(33, 9).Value = "sd"
$ws.Cells.Item(33, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(37, 9).Value = "sd"
$ws.Cells.Item(37, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(38, 9).Value = "aa"
$ws.Cells.Item(38, 10).Value = "Agree/Accept"
$ws.Cells.Item(39, 9).Value = "sd"
$ws.Cells.Item(39, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(46, 9).Value = "sd"
$ws.Cells.Item(46, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(47, 9).Value = "sd"
$ws.Cells.Item(47, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(48, 9).Value = "b"
$ws.Cells.Item(48, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(57, 9).Value = "sd"
$ws.Cells.Item(57, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(60, 9).Value = "sv"
$ws.Cells.Item(60, 10).Value = "Statement-opinion"
$ws.Cells.Item(65, 9).Value = "sd"
$ws.Cells.Item(65, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(67, 9).Value = "aa"
$ws.Cells.Item(67, 10).Value = "Agree/Accept"
$ws.Cells.Item(70, 9).Value = "sd"
$ws.Cells.Item(70, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(73, 9).Value = "b"
$ws.Cells.Item(73, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(83, 9).Value = "ba"
$ws.Cells.Item(83, 10).Value = "Appreciation"
$ws.Cells.Item(86, 9).Value = "b"
$ws.Cells.Item(86, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(95, 9).Value = "sd"
$ws.Cells.Item(95, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(115, 9).Value = "sd"
$ws.Cells.Item(115, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(122, 9).Value = "sd"
$ws.Cells.Item(122, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(130, 9).Value = "%"
$ws.Cells.Item(130, 10).Value = "Uninterpretable"
$ws.Cells.Item(131, 9).Value = "aa"
$ws.Cells.Item(131, 10).Value = "Agree/Accept"
$ws.Cells.Item(139, 9).Value = "%"
$ws.Cells.Item(139, 10).Value = "Uninterpretable"
$ws.Cells.Item(142, 9).Value = "aa"
$ws.Cells.Item(142, 10).Value = "Agree/Accept"
$ws.Cells.Item(156, 9).Value = "ba"
$ws.Cells.Item(156, 10).Value = "Appreciation"
$ws.Cells.Item(160, 9).Value = "b"
$ws.Cells.Item(160, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(162, 9).Value = "b"
$ws.Cells.Item(162, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(180, 9).Value = "aa"
$ws.Cells.Item(180, 10).Value = "Agree/Accept"
$ws.Cells.Item(185, 9).Value = "sd"
$ws.Cells.Item(185, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(187, 9).Value = "sd"
$ws.Cells.Item(187, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(209, 9).Value = "aa"
$ws.Cells.Item(209, 10).Value = "Agree/Accept"
$ws.Cells.Item(216, 9).Value = "sv"
$ws.Cells.Item(216, 10).Value = "Statement-opinion"
$ws.Cells.Item(224, 9).Value = "%"
$ws.Cells.Item(224, 10).Value = "Uninterpretable"
$ws.Cells.Item(226, 9).Value = "b"
$ws.Cells.Item(226, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(242, 9).Value = "sv"
$ws.Cells.Item(242, 10).Value = "Statement-opinion"
$ws.Cells.Item(264, 9).Value = "sv"
$ws.Cells.Item(264, 10).Value = "Statement-opinion"
$ws.Cells.Item(274, 9).Value = "b"
$ws.Cells.Item(274, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(276, 9).Value = "sv"
$ws.Cells.Item(276, 10).Value = "Statement-opinion"
$ws.Cells.Item(285, 9).Value = "sv"
$ws.Cells.Item(285, 10).Value = "Statement-opinion"
$ws.Cells.Item(292, 9).Value = "sv"
$ws.Cells.Item(292, 10).Value = "Statement-opinion"
$ws.Cells.Item(301, 9).Value = "sd"
$ws.Cells.Item(301, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(305, 9).Value = "sv"
$ws.Cells.Item(305, 10).Value = "Statement-opinion"
$ws.Cells.Item(331, 9).Value = "sd"
$ws.Cells.Item(331, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(336, 9).Value = "sv"
$ws.Cells.Item(336, 10).Value = "Statement-opinion"
$ws.Cells.Item(353, 9).Value = "sv"
$ws.Cells.Item(353, 10).Value = "Statement-opinion"
$ws.Cells.Item(365, 9).Value = "sv"
$ws.Cells.Item(365, 10).Value = "Statement-opinion"
$ws.Cells.Item(367, 9).Value = "sv"
$ws.Cells.Item(367, 10).Value = "Statement-opinion"
$ws.Cells.Item(373, 9).Value = "aa"
$ws.Cells.Item(373, 10).Value = "Agree/Accept"
$ws.Cells.Item(374, 9).Value = "sd"
$ws.Cells.Item(374, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(376, 9).Value = "aa"
$ws.Cells.Item(376, 10).Value = "Agree/Accept"
$ws.Cells.Item(382, 9).Value = "aa"
$ws.Cells.Item(382, 10).Value = "Agree/Accept"
$ws.Cells.Item(383, 9).Value = "sd"
$ws.Cells.Item(383, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(386, 9).Value = "sv"
$ws.Cells.Item(386, 10).Value = "Statement-opinion"
$ws.Cells.Item(417, 9).Value = "aa"
$ws.Cells.Item(417, 10).Value = "Agree/Accept"
$ws.Cells.Item(418, 9).Value = "sd"
$ws.Cells.Item(418, 10).Value = "Statement-non-opinion"
